# Generate Report for Handoff
# Adds a new localization-status row (for file
# c2307cba-fc6b-4884-b300-a3ef154fd862.md) to the Overview, zh-cn and
# de-de worksheets/tables.

$wb = $excel.ActiveWorkbook

$DATEFMT = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("A3").Value = "c2307cba-fc6b-4884-b300-a3ef154fd862.md"
$ws.Range("C3").Value = ".md"
# D3 (Publish URL) left blank on purpose - matches source data.
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-20 04:43:52"
$ws.Range("G3").NumberFormat = $DATEFMT

$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9efac57701012e2b0b9113da5e7db2bfd949a972/e2e/c2307cba-fc6b-4884-b300-a3ef154fd862.md"
$ws.Hyperlinks.Add($ws.Range("B3"), $url1, "", "", "e2e\c2307cba-fc6b-4884-b300-a3ef154fd862.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "c2307cba-fc6b-4884-b300-a3ef154fd862.d817dc4bc2fc4170176d0b9403a088abca196000.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-20 04:43:48"
$ws.Range("H3").NumberFormat = $DATEFMT
# I3, J3 (Latest Target File, Latest Handback File) left blank.
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = $DATEFMT
# L3 (Reference Tokens) left blank.
$ws.Range("M3").Value = "'True"
# N3 (Dependency From) left blank.
$ws.Range("O3").Value = "'False"
# P3 (Error Detail) left blank.

$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9efac57701012e2b0b9113da5e7db2bfd949a972/e2e/c2307cba-fc6b-4884-b300-a3ef154fd862.md"
$ws.Hyperlinks.Add($ws.Range("A3"), $url2, "", "", "c2307cba-fc6b-4884-b300-a3ef154fd862.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "c2307cba-fc6b-4884-b300-a3ef154fd862.d817dc4bc2fc4170176d0b9403a088abca196000.de-de.xlf"
$ws.Range("H3").Value = "2016-08-20 04:43:52"
$ws.Range("H3").NumberFormat = $DATEFMT
# I3, J3 (Latest Target File, Latest Handback File) left blank.
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = $DATEFMT
# L3 (Reference Tokens) left blank.
$ws.Range("M3").Value = "'True"
# N3 (Dependency From) left blank.
$ws.Range("O3").Value = "'False"
# P3 (Error Detail) left blank.

$url3 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9efac57701012e2b0b9113da5e7db2bfd949a972/e2e/c2307cba-fc6b-4884-b300-a3ef154fd862.md"
$ws.Hyperlinks.Add($ws.Range("A3"), $url3, "", "", "c2307cba-fc6b-4884-b300-a3ef154fd862.md") | Out-Null

$wb.Save()
